$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3576.923
$ws.Range("I64").Value = 3625
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 3625
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -3377
$ws.Range("N64").Value = -3496
$ws.Range("H67").Value = 3576.923
$ws.Range("I67").Value = 3625
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 3625
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -2767
$ws.Range("N67").Value = -4716
$ws.Range("H100").Value = 2071.4285
$ws.Range("I100").Value = 1955.5555
$ws.Range("J100").Value = 2280
$ws.Range("K100").Value = 1955.5555
$ws.Range("L100").Value = 2280
$ws.Range("M100").Value = -1414.5555
$ws.Range("N100").Value = -3362
$ws.Range("H103").Value = 100768.9
$ws.Range("I103").Value = 143512.72
$ws.Range("J103").Value = 1033.3334
$ws.Range("K103").Value = 430538.16
$ws.Range("L103").Value = 3100.0002
$ws.Range("M103").Value = -429952.16
$ws.Range("N103").Value = -4272.0002
$ws.Range("H137").Value = 2951.9688
$ws.Range("I137").Value = 2723.3333
$ws.Range("K137").Value = 8169.999899999999
$ws.Range("M137").Value = -5619.999899999999
$ws.Range("H138").Value = 1639.2742
$ws.Range("I138").Value = 1261.2894
$ws.Range("J138").Value = 2237.75
$ws.Range("K138").Value = 3783.8682
$ws.Range("L138").Value = 6713.25
$ws.Range("M138").Value = 1356.1318
$ws.Range("N138").Value = -16993.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 503797.12
$ws.Range("I32").Value = 553587.3
$ws.Range("K32").Value = 553587.3
$ws.Range("M32").Value = -553300.3
$ws.Range("H61").Value = 3146.0715
$ws.Range("I61").Value = 2318.7144
$ws.Range("J61").Value = 3973.4285
$ws.Range("K61").Value = 2318.7144
$ws.Range("L61").Value = 3973.4285
$ws.Range("M61").Value = -2106.7144
$ws.Range("N61").Value = -4397.4285
$ws.Range("H132").Value = 4221.6216
$ws.Range("I132").Value = 2964.4
$ws.Range("J132").Value = 5700.706
$ws.Range("K132").Value = 8893.200000000001
$ws.Range("L132").Value = 17102.118
$ws.Range("M132").Value = -6363.200000000001
$ws.Range("N132").Value = -22162.118
$ws.Range("H136").Value = 3146.0715
$ws.Range("I136").Value = 2318.7144
$ws.Range("J136").Value = 3973.4285
$ws.Range("K136").Value = 6956.1432
$ws.Range("L136").Value = 11920.2855
$ws.Range("M136").Value = -4406.1432
$ws.Range("N136").Value = -17020.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7355441
$ws.Range("I105").Value = 9618038
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 9618038
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -9616291
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1370
$ws.Range("J16").Value = 1660
$ws.Range("L16").Value = 1660
$ws.Range("N16").Value = -2234
$ws.Range("H31").Value = 6564.7393
$ws.Range("I31").Value = 1402.909
$ws.Range("J31").Value = 11296.417
$ws.Range("K31").Value = 1402.909
$ws.Range("L31").Value = 11296.417
$ws.Range("M31").Value = -1107.909
$ws.Range("N31").Value = -11886.417
$ws.Range("H34").Value = 6564.7393
$ws.Range("I34").Value = 1402.909
$ws.Range("J34").Value = 11296.417
$ws.Range("K34").Value = 1402.909
$ws.Range("L34").Value = 11296.417
$ws.Range("M34").Value = -1200.909
$ws.Range("N34").Value = -11700.417
$ws.Range("H105").Value = 1742.8572
$ws.Range("I105").Value = 1680
$ws.Range("J105").Value = 1900
$ws.Range("K105").Value = 1680
$ws.Range("L105").Value = 1900
$ws.Range("M105").Value = 67
$ws.Range("N105").Value = -5394
$ws.Range("H113").Value = 1370
$ws.Range("J113").Value = 1660
$ws.Range("L113").Value = 1660
$ws.Range("N113").Value = -6000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 615.4091
$ws.Range("J5").Value = 2240
$ws.Range("L5").Value = 6720
$ws.Range("N5").Value = -6944
$ws.Range("H12").Value = 76
$ws.Range("J12").Value = 94
$ws.Range("L12").Value = 282
$ws.Range("N12").Value = -628
$ws.Range("H22").Value = 1162.9025
$ws.Range("I22").Value = 937.9
$ws.Range("J22").Value = 1235.4839
$ws.Range("K22").Value = 2813.7
$ws.Range("L22").Value = 3706.4517
$ws.Range("M22").Value = -2644.7
$ws.Range("N22").Value = -4044.4517
$ws.Range("H27").Value = 1162.9025
$ws.Range("I27").Value = 937.9
$ws.Range("J27").Value = 1235.4839
$ws.Range("K27").Value = 2813.7
$ws.Range("L27").Value = 3706.4517
$ws.Range("M27").Value = -2711.7
$ws.Range("N27").Value = -3910.4517
$ws.Range("H34").Value = 12820979
$ws.Range("I34").Value = 208
$ws.Range("J34").Value = 14706387
$ws.Range("K34").Value = 624
$ws.Range("L34").Value = 44119161
$ws.Range("M34").Value = -540
$ws.Range("N34").Value = -44119329
$ws.Range("H35").Value = 3911.2
$ws.Range("J35").Value = 4814
$ws.Range("L35").Value = 14442
$ws.Range("N35").Value = -15018
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H58").Value = 2063.276
$ws.Range("I58").Value = 975
$ws.Range("J58").Value = 2102.1428
$ws.Range("K58").Value = 2925
$ws.Range("L58").Value = 6306.428400000001
$ws.Range("M58").Value = -2797
$ws.Range("N58").Value = -6562.428400000001
$ws.Range("H131").Value = 1027.6595
$ws.Range("J131").Value = 1103.0952
$ws.Range("L131").Value = 3309.2856
$ws.Range("N131").Value = -13389.2856
$ws.Range("H135").Value = 615.4091
$ws.Range("J135").Value = 2240
$ws.Range("L135").Value = 20160
$ws.Range("N135").Value = -25230
$ws.Range("H141").Value = 8900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 73086.14
$ws.Range("I113").Value = 101419
$ws.Range("J113").Value = 2254
$ws.Range("K113").Value = 101419
$ws.Range("L113").Value = 2254
$ws.Range("M113").Value = -99249
$ws.Range("N113").Value = -6594
$ws.Range("H122").Value = 2579.8
$ws.Range("I122").Value = 2249.5
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 6748.5
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -4298.5
$ws.Range("N122").Value = -13300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1871.3334
$ws.Range("I61").Value = 727.4167
$ws.Range("K61").Value = 727.4167
$ws.Range("M61").Value = -525.4167
$ws.Range("H74").Value = 32000
$ws.Range("J74").Value = 32000
$ws.Range("L74").Value = 32000
$ws.Range("N74").Value = -33996
$ws.Range("H77").Value = 32000
$ws.Range("J77").Value = 32000
$ws.Range("L77").Value = 96000
$ws.Range("N77").Value = -105984
$ws.Range("H113").Value = 1871.3334
$ws.Range("I113").Value = 727.4167
$ws.Range("K113").Value = 727.4167
$ws.Range("M113").Value = 1442.5833

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 70102
$ws.Range("J59").Value = 70102
$ws.Range("L59").Value = 70102
$ws.Range("N59").Value = -71578
$ws.Range("H113").Value = 1091.579
$ws.Range("I113").Value = 1262.2142
$ws.Range("J113").Value = 613.8
$ws.Range("K113").Value = 3786.6426
$ws.Range("L113").Value = 1841.4
$ws.Range("M113").Value = -1616.6426
$ws.Range("N113").Value = -6181.4
